# Documentation update: rename the "prevTaskBook" variable shown in the
# UndoRedoStartingStackDiagram's "DeleteCommand" state table to
# "prevOrganizer" (reflects the Addressbook/TaskBook -> PrioriTask model
# rename called out in the commit message).
#
# The text lives in a table cell on the diagram slide, in the second
# row's second paragraph ("prevTaskBook = s3"). We scan every table on
# every slide for the word and overwrite just the matching run's text so
# the rest of the cell ("targetIndex = 5" / " = s3") is left alone.

$p = $ppt.ActivePresentation

$oldWord = "prevTaskBook"
$newWord = "prevOrganizer"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $s.Shapes.Count; $shapeIdx++) {
        $sh = $s.Shapes.Item($shapeIdx)

        if ($sh.HasTable) {
            $tbl = $sh.Table

            for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
                for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                    $cellRange = $tbl.Cell($r, $c).Shape.TextFrame.TextRange
                    $cellText = $cellRange.Text
                    $idx = $cellText.IndexOf($oldWord)
                    if ($idx -ge 0) {
                        $cellRange.Characters($idx + 1, $oldWord.Length).Text = $newWord
                    }
                }
            }
        }
    }
}
